$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C4").Value = "B"
$ws.Range("C4").Select()
